# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> bound to the Notes Master ("Office Theme" colours)
#   ppt/theme/theme2.xml  -> bound to the Slide Master  ("Integral" colours)
# The authored edit swaps the two themes' contents: the Slide Master's theme
# becomes "Office Theme" and the Notes Master's theme becomes "Integral".
#
# The only theme surface the PowerPoint object model exposes for writing is
# the Slide Master's theme (Presentation.SlideMaster.Theme /
# Presentation.Designs(1)), so we drive the swap through that: push the
# "Office Theme" colour scheme (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink)
# onto it, one RGB value at a time, exactly as a user re-colouring the theme
# via Design > Variants > Colors would.

$p = $ppt.ActivePresentation

$design = $p.Designs.Item(1)
$theme = $design.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# Office Theme colour scheme, in ThemeColorScheme.Item(index) order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink.
# RGB values are the standard COM RGB() encoding (0x00BBGGRR) of the
# srgbClr hex values from the target theme XML.
$colors.Item(1).RGB  = 0         # dk1      000000
$colors.Item(2).RGB  = 16777215  # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388   # dk2      44546A
$colors.Item(4).RGB  = 15132391  # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939  # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501   # accent2  ED7D31
$colors.Item(7).RGB  = 10855845  # accent3  A5A5A5
$colors.Item(8).RGB  = 49407     # accent4  FFC000
$colors.Item(9).RGB  = 12874308  # accent5  4472C4
$colors.Item(10).RGB = 4697456   # accent6  70AD47
$colors.Item(11).RGB = 12673797  # hlink    0563C1
$colors.Item(12).RGB = 7491477   # folHlink 954F72
